# harmonized similar tags to be the same
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("isa_template")

# Row 12: ER protocol tag ("growth protocol" -> "growth")
$ws.Range("B12").Value = "growth"

# Row 13: Tags term -- EFO:0003789 -> GO:0040007, and the NCIT URL -> short-form NCIT:C14258
$ws.Range("B13").Value = "GO:0040007"
$ws.Range("D13").Value = "NCIT:C14258"

# Row 14: the now-unused Term Source REF ("EFO") and Term Accession Number ("NCIT") cells are cleared
$ws.Range("B14").ClearContents() | Out-Null
$ws.Range("D14").ClearContents() | Out-Null

# Row 13 no longer needs its taller custom row height now that the text is shorter
$ws.Rows.Item(13).AutoFit() | Out-Null

# Reflect where the author ended up with their selection after editing
$ws.Activate() | Out-Null
$ws.Range("E19").Select() | Out-Null
